# Adds a "Tipo_Unidade" column (Unidade/Kg) to the "Estoque" and
# "Removidos" sheets, inserted right after "Quantidade" and before
# "Preco". Also updates a few records on the "Estoque" sheet and
# appends two new stock rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Estoque"
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Estoque")

# Insert a new column D, shifting Preco/Categoria/Data_Criacao/
# Data_Alteracao one column to the right (D->E, E->F, F->G, G->H).
$ws.Columns.Item(4).Insert(-4161)

# New header
$ws.Range("D1").Value = "Tipo_Unidade"

# Default every existing item to "Unidade" ...
$ws.Range("D2:D15").Value = "Unidade"
# ... except the ones that are sold by weight.
$ws.Range("D4").Value = "Kg"

# A couple of records were edited along with this change.
$ws.Range("A2").Value = "Argamassa"
$ws.Range("H2").Value = "13/11/2025 23:42"
$ws.Range("H4").Value = "13/11/2025 23:39"

# Two new items were registered.
$ws.Range("A16").Value = "argam"
$ws.Range("B16").Value = "ID_15"
$ws.Range("C16").Value = 23
$ws.Range("D16").Value = "Kg"
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = "cimento"
$ws.Range("G16").Value = "13/11/2025 23:41"
$ws.Range("H16").Value = "13/11/2025 23:41"

$ws.Range("A17").Value = "gugu"
$ws.Range("B17").Value = "ID_16"
$ws.Range("C17").Value = 72
$ws.Range("D17").Value = "Kg"
$ws.Range("E17").Value = 30
$ws.Range("F17").Value = "bucha"
$ws.Range("G17").Value = "13/11/2025 23:43"
$ws.Range("H17").Value = "13/11/2025 23:43"

# ---------------------------------------------------------------------
# Sheet "Removidos"
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Removidos")

# Same column insertion as above (D->E, E->F, F->G, G->H, H->I).
$ws2.Columns.Item(4).Insert(-4161)

$ws2.Range("D1").Value = "Tipo_Unidade"
$ws2.Range("D2:D21").Value = "Unidade"

# The "Preco" column (now E) was blank for items that were removed
# before a price had ever been recorded; keep those cells empty
# (the column insert/shift otherwise leaves a stray empty string).
$ws2.Range("E2:E15").ClearContents()
$ws2.Range("E19").ClearContents()
